$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new product unit row
$ws.Range("A12").Value = "MILLIGRAMS"

# Move the active selection to J3 (matches the post-edit cursor position)
$ws.Range("J3").Select()
